$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 5130
$ws.Range("I2").Value = 13933
$ws.Range("J2").Value = 57860
$ws.Range("K2").Value = 269
$ws.Range("L2").Value = 15635
$ws.Range("M2").Value = 983
$ws.Range("N2").Value = 9987
$ws.Range("O2").Value = 32
$ws.Range("P2").Value = 254
$ws.Range("Q2").Value = 95
$ws.Range("R2").Value = 718
$ws.Range("S2").Value = 6130
$ws.Range("T2").Value = 10207
$ws.Range("U2").Value = 765
$ws.Range("V2").Value = 89185
$ws.Range("W2").Value = 28
$ws.Range("X2").Value = 89213
$ws.Range("Y2").Value = 138
$ws.Range("Z2").Value = 1369
$ws.Range("AA2").Value = 606
